$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bets")

# Copy the date format used in B2 onto the new date cells B3:B5
$ws.Range("B2").Copy()
$ws.Range("B3:B5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 45208
$ws.Range("C3").Value = 1
$ws.Range("D3").Formula = "=F2"
$ws.Range("E3").Value = 958
$ws.Range("F3").Formula = "=D3+E3"

# Row 4
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 45208
$ws.Range("C4").Value = 0
$ws.Range("D4").Formula = "=F3"
$ws.Range("E4").Value = -485
$ws.Range("F4").Formula = "=D4+E4"

# Row 5
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 45209
$ws.Range("C5").Value = 1
$ws.Range("D5").Formula = "=F4"
$ws.Range("E5").Value = 382
$ws.Range("F5").Formula = "=D5+E5"

$ws.Range("G5").Select()
